$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: safe to assign directly
$plainValues = @{
    'D2' = '23.864.86'
    'E2' = '  -3.26%  '
    'D3' = '1.620.97'
    'E3' = '  -3.26%  '
    'E4' = '  -0.17%  '
    'E5' = '  -1.84%  '
    'E6' = '  -0.12%  '
    'E7' = '  -0.19%  '
    'E8' = '  -2.48%  '
    'B9' = 'BinanceUSD'
    'C9' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'E9' = '  -0.15%  '
    'B10' = 'Polygon'
    'C10' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'E10' = '  -2.37%  '
    'E11' = '  -2.53%  '
    'E12' = '  -2.24%  '
    'E13' = '  -4.89%  '
    'E15' = '  -2.07%  '
    'E16' = '  -2.58%  '
    'D17' = '1.620.62'
    'E17' = '  -3.49%  '
    'E18' = '  -0.18%  '
    'E19' = '  -1.16%  '
    'E20' = '  -5.09%  '
    'E21' = '  -3.62%  '
    'E22' = '  -0.31%  '
    'E23' = '  -3.72%  '
    'D24' = '23.852.03'
    'E24' = '  -3.33%  '
    'E25' = '  +4.00%  '
    'E26' = '  +3.56%  '
    'E27' = '  -3.57%  '
    'E28' = '  -2.09%  '
    'E29' = '  -4.38%  '
    'E30' = '  -10.16%  '
    'E31' = '  -5.73%  '
    'E32' = '  -0.16%  '
    'D33' = '1.801.29'
    'E33' = '  -3.33%  '
    'E34' = '  -1.80%  '
    'E35' = '  -0.65%  '
    'E36' = '  -6.42%  '
    'E37' = '  -5.56%  '
    'E38' = '  -4.55%  '
    'E39' = '  -4.91%  '
    'E40' = '  +0.66%  '
    'E41' = '  +0.50%  '
    'E42' = '  -6.47%  '
    'E43' = '  -5.03%  '
    'E44' = '  -3.22%  '
    'E45' = '  -2.77%  '
    'E46' = '  -3.50%  '
    'E47' = '  -2.38%  '
    'E48' = '  -0.11%  '
    'E49' = '  -4.87%  '
    'E50' = '  -3.03%  '
    'E51' = '  -8.17%  '
}
foreach ($addr in $plainValues.Keys) {
    $ws.Range($addr).Value = $plainValues[$addr]
}

# Numeric-looking text values that Excel would otherwise auto-convert to numbers:
# force Text number format while assigning, then restore default style so no
# stray style attribute is left on the cell.
$textForcedValues = @{
    'D4' = '0.9990'
    'D5' = '307.85'
    'D6' = '1.000'
    'D7' = '0.3930'
    'D8' = '0.3847'
    'D9' = '0.9993'
    'D10' = '1.370'
    'D11' = '49.60'
    'D12' = '0.08455'
    'D13' = '23.98'
    'D14' = '7.053'
    'D15' = '7.555'
    'D16' = '0.00001281'
    'D18' = '93.80'
    'D19' = '0.06932'
    'D20' = '20.09'
    'D21' = '6.824'
    'D22' = '0.9982'
    'D23' = '13.42'
    'D25' = '2.435'
    'D26' = '2.883'
    'D27' = '22.25'
    'D28' = '156.51'
    'D29' = '139.92'
    'D30' = '5.270'
    'D31' = '7.864'
    'D32' = '2.492'
    'D34' = '0.08141'
    'D35' = '0.9847'
    'D36' = '0.02889'
    'D37' = '6.589'
    'D38' = '0.2679'
    'D39' = '0.09152'
    'D40' = '10.36'
    'D41' = '13.57'
    'D42' = '1.425'
    'D43' = '0.7513'
    'D44' = '16.05'
    'D45' = '0.6905'
    'D46' = '2.476'
    'D47' = '4.075'
    'D48' = '0.9996'
    'D49' = '0.08249'
    'D50' = '133.66'
    'D51' = '1.220'
}
foreach ($addr in $textForcedValues.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedValues[$addr]
    $cell.Style = "Normal"
}
